$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-06-24 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-25 Tuesday", 2)

# Update the division problems in the table, identified by fixed cell position
# so replacements don't collide with each other (some new values equal other
# cells' old values).
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "38÷3=12, 2"
$t.Cell(1, 2).Range.Text = "29÷9=3, 2"
$t.Cell(1, 3).Range.Text = "56÷6=9, 2"
$t.Cell(1, 4).Range.Text = "66÷7=9, 3"
$t.Cell(1, 5).Range.Text = "67÷5=13, 2"

$t.Cell(5, 1).Range.Text = "48÷5=9, 3"
$t.Cell(5, 2).Range.Text = "48÷4=12, 0"
$t.Cell(5, 3).Range.Text = "67÷2=33, 1"
$t.Cell(5, 4).Range.Text = "99÷9=11, 0"
$t.Cell(5, 5).Range.Text = "90÷7=12, 6"

$t.Cell(9, 1).Range.Text = "40÷9=4, 4"
$t.Cell(9, 2).Range.Text = "63÷9=7, 0"
$t.Cell(9, 3).Range.Text = "93÷3=31, 0"
$t.Cell(9, 4).Range.Text = "96÷5=19, 1"
$t.Cell(9, 5).Range.Text = "55÷5=11, 0"

$t.Cell(13, 1).Range.Text = "19÷2=9, 1"
$t.Cell(13, 2).Range.Text = "26÷6=4, 2"
$t.Cell(13, 3).Range.Text = "24÷8=3, 0"
$t.Cell(13, 4).Range.Text = "44÷5=8, 4"
$t.Cell(13, 5).Range.Text = "49÷7=7, 0"

$t.Cell(17, 1).Range.Text = "86÷6=14, 2"
$t.Cell(17, 2).Range.Text = "96÷6=16, 0"
$t.Cell(17, 3).Range.Text = "56÷3=18, 2"
$t.Cell(17, 4).Range.Text = "82÷7=11, 5"
$t.Cell(17, 5).Range.Text = "17÷6=2, 5"
